$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Resumen")
$ws2 = $wb.Worksheets.Item("Solucion")
$ws3 = $wb.Worksheets.Item("Metricas")

# --- Resumen sheet updates ---
$ws1.Range("B2").Value = "Z2"
$ws1.Range("C2").Value = 555.7045410314629

# --- Solucion sheet updates (Pedido / Salida reassignments) ---
$ws2.Range("A2").Value = "Pedido_23"
$ws2.Range("A3").Value = "Pedido_49"
$ws2.Range("A4").Value = "Pedido_10"
$ws2.Range("A5").Value = "Pedido_8"
$ws2.Range("B5").Value = "S011"
$ws2.Range("A6").Value = "Pedido_20"
$ws2.Range("B6").Value = "S031"
$ws2.Range("A7").Value = "Pedido_36"
$ws2.Range("B7").Value = "S051"
$ws2.Range("A8").Value = "Pedido_41"
$ws2.Range("B8").Value = "S042"
$ws2.Range("A9").Value = "Pedido_32"
$ws2.Range("B9").Value = "S002"
$ws2.Range("A10").Value = "Pedido_26"
$ws2.Range("A11").Value = "Pedido_38"
$ws2.Range("B11").Value = "S052"
$ws2.Range("A12").Value = "Pedido_37"
$ws2.Range("A13").Value = "Pedido_5"
$ws2.Range("B13").Value = "S032"
$ws2.Range("A14").Value = "Pedido_9"
$ws2.Range("B14").Value = "S043"
$ws2.Range("A15").Value = "Pedido_45"
$ws2.Range("B15").Value = "S003"
$ws2.Range("A16").Value = "Pedido_21"
$ws2.Range("B16").Value = "S023"
$ws2.Range("A17").Value = "Pedido_17"
$ws2.Range("B17").Value = "S053"
$ws2.Range("A18").Value = "Pedido_28"
$ws2.Range("A19").Value = "Pedido_39"
$ws2.Range("B19").Value = "S013"
$ws2.Range("A20").Value = "Pedido_55"
$ws2.Range("B20").Value = "S044"
$ws2.Range("A21").Value = "Pedido_24"
$ws2.Range("B21").Value = "S004"
$ws2.Range("A22").Value = "Pedido_27"
$ws2.Range("B22").Value = "S024"
$ws2.Range("A23").Value = "Pedido_6"
$ws2.Range("B23").Value = "S054"
$ws2.Range("A24").Value = "Pedido_4"
$ws2.Range("B24").Value = "S014"
$ws2.Range("A25").Value = "Pedido_31"
$ws2.Range("B25").Value = "S045"
$ws2.Range("A26").Value = "Pedido_40"
$ws2.Range("A27").Value = "Pedido_54"
$ws2.Range("B27").Value = "S055"
$ws2.Range("A28").Value = "Pedido_2"
$ws2.Range("A29").Value = "Pedido_58"
$ws2.Range("B29").Value = "S005"
$ws2.Range("A30").Value = "Pedido_25"
$ws2.Range("B30").Value = "S015"
$ws2.Range("A31").Value = "Pedido_43"
$ws2.Range("B31").Value = "S046"
$ws2.Range("A32").Value = "Pedido_56"
$ws2.Range("B32").Value = "S035"
$ws2.Range("A33").Value = "Pedido_16"
$ws2.Range("B33").Value = "S056"
$ws2.Range("A34").Value = "Pedido_33"
$ws2.Range("B34").Value = "S006"
$ws2.Range("A35").Value = "Pedido_46"
$ws2.Range("B35").Value = "S026"
$ws2.Range("A36").Value = "Pedido_59"
$ws2.Range("B36").Value = "S036"
$ws2.Range("A37").Value = "Pedido_42"
$ws2.Range("B37").Value = "S016"
$ws2.Range("A38").Value = "Pedido_48"
$ws2.Range("B38").Value = "S047"
$ws2.Range("A39").Value = "Pedido_13"
$ws2.Range("A40").Value = "Pedido_35"
$ws2.Range("B40").Value = "S007"
$ws2.Range("B41").Value = "S037"
$ws2.Range("A42").Value = "Pedido_60"
$ws2.Range("B42").Value = "S057"
$ws2.Range("A43").Value = "Pedido_7"
$ws2.Range("B43").Value = "S028"
$ws2.Range("A44").Value = "Pedido_57"
$ws2.Range("B44").Value = "S017"
$ws2.Range("A45").Value = "Pedido_12"
$ws2.Range("B45").Value = "S048"
$ws2.Range("A46").Value = "Pedido_47"
$ws2.Range("B46").Value = "S058"
$ws2.Range("A47").Value = "Pedido_1"
$ws2.Range("B47").Value = "S008"
$ws2.Range("A48").Value = "Pedido_19"
$ws2.Range("B48").Value = "S038"
$ws2.Range("A49").Value = "Pedido_52"
$ws2.Range("B49").Value = "S049"
$ws2.Range("A50").Value = "Pedido_18"
$ws2.Range("B50").Value = "S018"
$ws2.Range("A51").Value = "Pedido_34"
$ws2.Range("B51").Value = "S029"
$ws2.Range("A52").Value = "Pedido_30"
$ws2.Range("B52").Value = "S059"
$ws2.Range("A53").Value = "Pedido_22"
$ws2.Range("B53").Value = "S039"
$ws2.Range("A54").Value = "Pedido_53"
$ws2.Range("B54").Value = "S009"
$ws2.Range("A55").Value = "Pedido_50"
$ws2.Range("B55").Value = "S019"
$ws2.Range("A56").Value = "Pedido_29"
$ws2.Range("B56").Value = "S050"
$ws2.Range("A57").Value = "Pedido_44"
$ws2.Range("B57").Value = "S010"
$ws2.Range("A58").Value = "Pedido_11"
$ws2.Range("B58").Value = "S030"
$ws2.Range("A59").Value = "Pedido_3"
$ws2.Range("B59").Value = "S020"
$ws2.Range("A60").Value = "Pedido_15"
$ws2.Range("A61").Value = "Pedido_14"

# --- Metricas sheet updates (Tiempo values) ---
$ws3.Range("B2").Value = 555.7045410314629
$ws3.Range("B3").Value = 555.7045410314629
$ws3.Range("B4").Value = 549.2799654016652
